$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("A2").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A4").Value = 3
$ws.Range("B2").Copy()
$ws.Range("B4").PasteSpecial(-4122)
$ws.Range("B4").Value = '自我卡通貼紙'
$ws.Range("B2").Copy()
$ws.Range("C4").PasteSpecial(-4122)
$ws.Range("C4").Value = '自我卡通貼紙'
$ws.Range("B2").Copy()
$ws.Range("D4").PasteSpecial(-4122)
$ws.Range("D4").Value = '00003.png'
$ws.Range("B2").Copy()
$ws.Range("E4").PasteSpecial(-4122)
$ws.Range("E4").Value = 'young2.jpg'
$ws.Range("B2").Copy()
$ws.Range("H4").PasteSpecial(-4122)
$ws.Range("H4").Value = 'Making a playful peace sign with both hands and winking. Tearful eyes and slightly trembling lips, showing a cute crying expression. Arms wide open in a warm, enthusiastic hug pose. Lying on their side asleep, resting on a tiny pillow with a sweet smile. Pointing forward with confidence, surrounded by shining visual effects. Blowing a kiss, with heart symbols floating around. Maintain the chibi aesthetic. Exaggerated, expressive big eyes. Soft facial lines. Background: Vibrant red with star or colorful confetti elements for decoration. Leave some clean white space around each sticker. Aspect ratio: 9:16'

# Row 5
$ws.Range("A2").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A5").Value = 4
$ws.Range("B2").Copy()
$ws.Range("B5").PasteSpecial(-4122)
$ws.Range("B5").Value = '更換衣服'
$ws.Range("B2").Copy()
$ws.Range("C5").PasteSpecial(-4122)
$ws.Range("C5").Value = '更換衣服'
$ws.Range("B2").Copy()
$ws.Range("D5").PasteSpecial(-4122)
$ws.Range("D5").Value = '00004.png'
$ws.Range("B2").Copy()
$ws.Range("E5").PasteSpecial(-4122)
$ws.Range("E5").Value = 'self01.jpg'
$ws.Range("B2").Copy()
$ws.Range("H5").PasteSpecial(-4122)
$ws.Range("H5").Value = '图中人物，更換上另一图的衣服。'

# Row 6
$ws.Range("A2").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("A6").Value = 5
$ws.Range("B2").Copy()
$ws.Range("B6").PasteSpecial(-4122)
$ws.Range("B6").Value = '3D figurine'
$ws.Range("B2").Copy()
$ws.Range("C6").PasteSpecial(-4122)
$ws.Range("C6").Value = '3D figurine'
$ws.Range("B2").Copy()
$ws.Range("D6").PasteSpecial(-4122)
$ws.Range("D6").Value = '00005.png'
$ws.Range("B2").Copy()
$ws.Range("E6").PasteSpecial(-4122)
$ws.Range("E6").Value = 'self01.jpg'
$ws.Range("B2").Copy()
$ws.Range("H6").PasteSpecial(-4122)
$ws.Range("H6").Value = 'Create a 1/7 scale commercialized figurine of the characters in the picture, in a realistic style, in a real environment. The figurine is placed on a computer desk. The figurine has a round transparent acrylic base, with no text on the base. The content on the computer screen is a 3D modeling process of this figurine. Next to the computer screen is a toy packaging box, designed in a style reminiscent of high-quality collectible figures, printed with original artwork. The packaging features two-dimensional flat illustrations'

# Row 7
$ws.Range("A2").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("A7").Value = 6
$ws.Range("B2").Copy()
$ws.Range("B7").PasteSpecial(-4122)
$ws.Range("B7").Value = '向國家致敬'
$ws.Range("B2").Copy()
$ws.Range("C7").PasteSpecial(-4122)
$ws.Range("C7").Value = '向國家致敬'
$ws.Range("B2").Copy()
$ws.Range("D7").PasteSpecial(-4122)
$ws.Range("D7").Value = '00006.png'
$ws.Range("B2").Copy()
$ws.Range("E7").PasteSpecial(-4122)
$ws.Range("E7").Value = 'self01.jpg'
$ws.Range("B2").Copy()
$ws.Range("F7").PasteSpecial(-4122)
$ws.Range("F7").Value = 'flag.jpg'
$ws.Range("B2").Copy()
$ws.Range("H7").PasteSpecial(-4122)
$ws.Range("H7").Value = '图中人向國家昂首致敬，衷心的，有氣勢的。衣服，背景請配合國旗(見图)。註，人物樣子不能改，但可配合情景'

$excel.CutCopyMode = 0